# 6.2.1.1.xlsx update: add a new "2023" data column (P) to the single
# worksheet, mirroring the formatting already used for the "2022" column
# (O), and mark two previously-empty cells in the "total" row with a
# dash ("-") placeholder, right-aligned like the rest of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: header (year) + per-region values for 2023 ----------
# Copying the O-column cell first carries over its number format / font /
# border (xf) so the new column matches the existing year columns, then
# we overwrite just the value.

$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("P4").Value = 2023

$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("P5").Value = 48.2

$ws.Range("O6").Copy($ws.Range("P6"))
$ws.Range("P6").Value = 8.6767564891727478

$ws.Range("O7").Copy($ws.Range("P7"))
$ws.Range("P7").Value = 12.226605469730881

$ws.Range("O8").Copy($ws.Range("P8"))
$ws.Range("P8").Value = 78.520866131691164

$ws.Range("O9").Copy($ws.Range("P9"))
$ws.Range("P9").Value = 59.466452648968115

$ws.Range("O10").Copy($ws.Range("P10"))
$ws.Range("P10").Value = 26.635270208942913

$ws.Range("O11").Copy($ws.Range("P11"))
$ws.Range("P11").Value = 8.166450559693871

$ws.Range("O12").Copy($ws.Range("P12"))
$ws.Range("P12").Value = 74.601894583630667

$ws.Range("O13").Copy($ws.Range("P13"))
$ws.Range("P13").Value = 99.168063426054971

$ws.Range("O14").Copy($ws.Range("P14"))
$ws.Range("P14").Value = 70.956108992253434

# Row 3 is the thin header-separator row: it only needs the matching
# (empty) border style carried into the new column, no value.
$ws.Range("O3").Copy($ws.Range("P3"))
$ws.Range("P3").ClearContents()

# --- Fill the two previously blank "total" cells with a dash -----------
$ws.Range("D14").Value = "-"
$ws.Range("D14").HorizontalAlignment = -4152   # xlRight
$ws.Range("E14").Value = "-"
$ws.Range("E14").HorizontalAlignment = -4152   # xlRight

# --- Row-height touch-ups that came along with the new column ----------
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 13.5

Write-Output "Added 2023 column (P) and dash placeholders in D14:E14"
